$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data block (row 288), pushing the
# existing rows 288-368 down to 290-370.
$ws.Rows("288:289").Insert()

# Fixed columns shared by every data row in this sheet.
$mercado = "Agrícola del Norte S.A. de Arica"
$region = "Arica y Parinacota"
$origen = "Región de Arica y Parinacota"
$clasificacion = "Hortaliza"

# New row 288: "Primera" quality for the latest reporting date.
$ws.Cells.Item(288, 1).Value = 1
$ws.Cells.Item(288, 2).Value = $mercado
$ws.Cells.Item(288, 3).Value = $region
$ws.Cells.Item(288, 4).Value = 44736
$ws.Cells.Item(288, 5).Value = 15
$ws.Cells.Item(288, 6).Value = 100112032
$ws.Cells.Item(288, 7).Value = "Zapallo italiano"
$ws.Cells.Item(288, 8).Value = "Huracán"
$ws.Cells.Item(288, 9).Value = "Primera"
$ws.Cells.Item(288, 10).Value = 130
$ws.Cells.Item(288, 11).Value = 8000
$ws.Cells.Item(288, 12).Value = 9000
$ws.Cells.Item(288, 13).Value = 8500
$ws.Cells.Item(288, 14).Value = '$/caja 70 unidades'
$ws.Cells.Item(288, 15).Value = $origen
$ws.Cells.Item(288, 16).Value = 121
$ws.Cells.Item(288, 17).Value = 70
$ws.Cells.Item(288, 18).Value = $clasificacion

# New row 289: "Segunda" quality for the latest reporting date.
$ws.Cells.Item(289, 1).Value = 1
$ws.Cells.Item(289, 2).Value = $mercado
$ws.Cells.Item(289, 3).Value = $region
$ws.Cells.Item(289, 4).Value = 44736
$ws.Cells.Item(289, 5).Value = 15
$ws.Cells.Item(289, 6).Value = 100112032
$ws.Cells.Item(289, 7).Value = "Zapallo italiano"
$ws.Cells.Item(289, 8).Value = "Huracán"
$ws.Cells.Item(289, 9).Value = "Segunda"
$ws.Cells.Item(289, 10).Value = 150
$ws.Cells.Item(289, 11).Value = 6000
$ws.Cells.Item(289, 12).Value = 7000
$ws.Cells.Item(289, 13).Value = 6500
$ws.Cells.Item(289, 14).Value = '$/caja 100 unidades'
$ws.Cells.Item(289, 15).Value = $origen
$ws.Cells.Item(289, 16).Value = 65
$ws.Cells.Item(289, 17).Value = 100
$ws.Cells.Item(289, 18).Value = $clasificacion
